$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1236.425
$ws.Range("J17").Value = 1271.7368
$ws.Range("L17").Value = 3815.2104
$ws.Range("N17").Value = -4151.2104

$ws.Range("H55").Value = 233.91667
$ws.Range("J55").Value = 254.5
$ws.Range("L55").Value = 254.5
$ws.Range("N55").Value = -682.5

$ws.Range("H108").Value = 74500
$ws.Range("J108").Value = 74500
$ws.Range("L108").Value = 74500
$ws.Range("N108").Value = -82180

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H138").Value = 6098.528
$ws.Range("J138").Value = 5527.983
$ws.Range("L138").Value = 16583.949
$ws.Range("N138").Value = -26863.949

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2142.5557
$ws.Range("I74").Value = 2248.1875
$ws.Range("J74").Value = 1297.5
$ws.Range("K74").Value = 2248.1875
$ws.Range("L74").Value = 1297.5
$ws.Range("M74").Value = -1374.1875
$ws.Range("N74").Value = -3045.5

$ws.Range("H77").Value = 2142.5557
$ws.Range("I77").Value = 2248.1875
$ws.Range("J77").Value = 1297.5
$ws.Range("K77").Value = 11240.9375
$ws.Range("L77").Value = 6487.5
$ws.Range("M77").Value = -6872.9375
$ws.Range("N77").Value = -15223.5

$ws.Range("H88").Value = 3745.2727
$ws.Range("I88").Value = 3039.2
$ws.Range("K88").Value = 3039.2
$ws.Range("M88").Value = -2633.2

$ws.Range("H91").Value = 3745.2727
$ws.Range("I91").Value = 3039.2
$ws.Range("K91").Value = 3039.2
$ws.Range("M91").Value = -1635.2

$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("M108").ClearContents()

$ws.Range("H122").Value = 2762.5
$ws.Range("I122").Value = 1650
$ws.Range("J122").Value = 3133.3333
$ws.Range("K122").Value = 4950
$ws.Range("L122").Value = 9399.999899999999
$ws.Range("M122").Value = -2500
$ws.Range("N122").Value = -14299.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 29496.666
$ws.Range("J76").Value = 29496.666
$ws.Range("L76").Value = 29496.666
$ws.Range("N76").Value = -30126.666

$ws.Range("H79").Value = 29496.666
$ws.Range("J79").Value = 29496.666
$ws.Range("L79").Value = 29496.666
$ws.Range("N79").Value = -31680.666

$ws.Range("H107").Value = 5310.769
$ws.Range("I107").Value = 2380
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 2380
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = -460
$ws.Range("N107").Value = -13840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 920.3
$ws.Range("J107").Value = 1508.3846
$ws.Range("L107").Value = 1508.3846
$ws.Range("N107").Value = -5348.3846

$ws.Range("H132").Value = 1684.0834
$ws.Range("I132").Value = 1636.0883
$ws.Range("K132").Value = 4908.2649
$ws.Range("M132").Value = -2378.2649

$ws.Range("H134").Value = 4281.364
$ws.Range("I134").Value = 3709.5
$ws.Range("K134").Value = 11128.5
$ws.Range("M134").Value = -8593.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 833.7646999999999
$ws.Range("I5").Value = 605.6923
$ws.Range("J5").Value = 1575
$ws.Range("K5").Value = 1817.0769
$ws.Range("L5").Value = 4725
$ws.Range("M5").Value = -1705.0769
$ws.Range("N5").Value = -4949

$ws.Range("H7").Value = 249.28572
$ws.Range("I7").Value = 287.83334
$ws.Range("J7").Value = 18
$ws.Range("K7").Value = 863.5000200000001
$ws.Range("L7").Value = 54
$ws.Range("M7").Value = -751.5000200000001
$ws.Range("N7").Value = -278

$ws.Range("H40").Value = 172.71428
$ws.Range("I40").Value = 53
$ws.Range("K40").Value = 212
$ws.Range("M40").Value = -143

$ws.Range("H52").Value = 1471.8
$ws.Range("J52").Value = 1471.8
$ws.Range("L52").Value = 4415.4
$ws.Range("N52").Value = -4947.4

$ws.Range("H75").Value = 6833
$ws.Range("J75").Value = 10000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31996

$ws.Range("H78").Value = 6833
$ws.Range("J78").Value = 10000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99984

$ws.Range("H101").Value = 12500
$ws.Range("J101").Value = 12500
$ws.Range("L101").Value = 37500
$ws.Range("N101").Value = -42368

$ws.Range("H132").Value = 2308.9333
$ws.Range("J132").Value = 2777.6667
$ws.Range("L132").Value = 24999.0003
$ws.Range("N132").Value = -30059.0003

$ws.Range("H135").Value = 833.7646999999999
$ws.Range("I135").Value = 605.6923
$ws.Range("J135").Value = 1575
$ws.Range("K135").Value = 5451.2307
$ws.Range("L135").Value = 14175
$ws.Range("M135").Value = -2916.2307
$ws.Range("N135").Value = -19245

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 47000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 47000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 47000
$ws.Range("N62").Value = -48372
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 47000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 47000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 141000
$ws.Range("N65").Value = -147864
$ws.Range("M65").ClearContents()

$ws.Range("H132").Value = 3448.7273
$ws.Range("I132").Value = 3448.7273
$ws.Range("K132").Value = 10346.1819
$ws.Range("M132").Value = -7816.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 30000
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 30000
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H74").Value = 56666.668
$ws.Range("I74").Value = 50000
$ws.Range("K74").Value = 50000
$ws.Range("M74").Value = -49002

$ws.Range("H77").Value = 56666.668
$ws.Range("I77").Value = 50000
$ws.Range("K77").Value = 150000
$ws.Range("M77").Value = -145008

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1852.9375
$ws.Range("I81").Value = 1331.9286
$ws.Range("K81").Value = 2663.8572
$ws.Range("M81").Value = -1602.8572

$ws.Range("H84").Value = 1852.9375
$ws.Range("I84").Value = 1331.9286
$ws.Range("K84").Value = 13319.286
$ws.Range("M84").Value = -8015.286

$ws.Range("H110").Value = 133699.8
$ws.Range("J110").Value = 133699.8
$ws.Range("L110").Value = 133699.8
$ws.Range("N110").Value = -141879.8

$ws.Range("H132").Value = 6543.3794
$ws.Range("I132").Value = 6065.1113
$ws.Range("K132").Value = 18195.3339
$ws.Range("M132").Value = -15665.3339
